$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 1132.25
$ws.Range("I42").Value = 10.333333
$ws.Range("J42").Value = 1805.4
$ws.Range("K42").Value = 30.999999
$ws.Range("L42").Value = 5416.200000000001
$ws.Range("M42").Value = 199.000001
$ws.Range("N42").Value = -5876.200000000001
$ws.Range("H76").Value = 3599.3333
$ws.Range("I76").Value = 3527.8572
$ws.Range("K76").Value = 3527.8572
$ws.Range("M76").Value = -3212.8572
$ws.Range("H79").Value = 3599.3333
$ws.Range("I79").Value = 3527.8572
$ws.Range("K79").Value = 3527.8572
$ws.Range("M79").Value = -2435.8572
$ws.Range("H121").Value = 3622.818
$ws.Range("J121").Value = 3622.818
$ws.Range("L121").Value = 10868.454
$ws.Range("N121").Value = -14362.454
$ws.Range("H125").Value = 3624.5
$ws.Range("J125").Value = 4499.6665
$ws.Range("L125").Value = 40496.9985
$ws.Range("N125").Value = -45416.9985
$ws.Range("H131").Value = 4023.7058
$ws.Range("I131").Value = 700.9167
$ws.Range("K131").Value = 2102.7501
$ws.Range("M131").Value = 2937.2499
$ws.Range("H132").Value = 126026.26
$ws.Range("I132").Value = 362332.3
$ws.Range("K132").Value = 1086996.9
$ws.Range("M132").Value = -1084466.9
$ws.Range("H138").Value = 5264.5903
$ws.Range("I138").Value = 1228.5238
$ws.Range("K138").Value = 3685.5714
$ws.Range("M138").Value = 1454.4286
$ws.Range("H141").Value = 6973.4375
$ws.Range("I141").Value = 5827.143
$ws.Range("J141").Value = 14997.5
$ws.Range("K141").Value = 17481.429
$ws.Range("L141").Value = 44992.5
$ws.Range("M141").Value = -12301.429
$ws.Range("N141").Value = -55352.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3665.468
$ws.Range("I32").Value = 1719.6
$ws.Range("J32").Value = 14784.714
$ws.Range("K32").Value = 1719.6
$ws.Range("L32").Value = 14784.714
$ws.Range("M32").Value = -1432.6
$ws.Range("N32").Value = -15358.714
$ws.Range("H61").Value = 21403.375
$ws.Range("I61").Value = 21534.666
$ws.Range("J61").Value = 21324.6
$ws.Range("K61").Value = 21534.666
$ws.Range("L61").Value = 21324.6
$ws.Range("M61").Value = -21322.666
$ws.Range("N61").Value = -21748.6
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H74").Value = 1033.6897
$ws.Range("I74").Value = 399.66666
$ws.Range("J74").Value = 1481.2354
$ws.Range("K74").Value = 399.66666
$ws.Range("L74").Value = 1481.2354
$ws.Range("M74").Value = 474.33334
$ws.Range("N74").Value = -3229.2354
$ws.Range("H77").Value = 1033.6897
$ws.Range("I77").Value = 399.66666
$ws.Range("J77").Value = 1481.2354
$ws.Range("K77").Value = 1998.3333
$ws.Range("L77").Value = 7406.177
$ws.Range("M77").Value = 2369.6667
$ws.Range("N77").Value = -16142.177
$ws.Range("H88").Value = 168000
$ws.Range("I88").Value = 2000
$ws.Range("K88").Value = 2000
$ws.Range("M88").Value = -1594
$ws.Range("H91").Value = 168000
$ws.Range("I91").Value = 2000
$ws.Range("K91").Value = 2000
$ws.Range("M91").Value = -596
$ws.Range("H122").Value = 3613.6924
$ws.Range("I122").Value = 2348.7896
$ws.Range("K122").Value = 7046.3688
$ws.Range("M122").Value = -4596.3688
$ws.Range("H132").Value = 28177.64
$ws.Range("I132").Value = 34696.65
$ws.Range("K132").Value = 104089.95
$ws.Range("M132").Value = -101559.95
$ws.Range("H136").Value = 21403.375
$ws.Range("I136").Value = 21534.666
$ws.Range("J136").Value = 21324.6
$ws.Range("K136").Value = 64603.99800000001
$ws.Range("L136").Value = 63973.8
$ws.Range("M136").Value = -62053.99800000001
$ws.Range("N136").Value = -69073.79999999999
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("H80").Value = 474.66666
$ws.Range("J80").Value = 388.66666
$ws.Range("L80").Value = 388.66666
$ws.Range("N80").Value = -2384.66666
$ws.Range("H83").Value = 474.66666
$ws.Range("J83").Value = 388.66666
$ws.Range("L83").Value = 1943.3333
$ws.Range("N83").Value = -11927.3333
$ws.Range("H93").Value = 96542.28999999999
$ws.Range("J93").Value = 96542.28999999999
$ws.Range("L93").Value = 96542.28999999999
$ws.Range("N93").Value = -100286.29
$ws.Range("H118").Value = 54998
$ws.Range("J118").Value = 54998
$ws.Range("L118").Value = 54998
$ws.Range("N118").Value = -58312
$ws.Range("H134").Value = 2745.054
$ws.Range("I134").Value = 2255.7334
$ws.Range("K134").Value = 6767.2002
$ws.Range("M134").Value = -4232.2002
$ws.Range("N53").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("H60").Value = 175874.75
$ws.Range("I60").Value = 500
$ws.Range("J60").Value = 234333
$ws.Range("K60").Value = 500
$ws.Range("L60").Value = 234333
$ws.Range("M60").Value = 11
$ws.Range("N60").Value = -235355
$ws.Range("H86").Value = 5696.6
$ws.Range("I86").Value = 5615
$ws.Range("K86").Value = 5615
$ws.Range("M86").Value = -4492
$ws.Range("H89").Value = 5696.6
$ws.Range("I89").Value = 5615
$ws.Range("K89").Value = 28075
$ws.Range("M89").Value = -22459
$ws.Range("H99").Value = 11806.8125
$ws.Range("I99").Value = 20318.5
$ws.Range("K99").Value = 20318.5
$ws.Range("M99").Value = -18820.5
$ws.Range("H126").Value = 11806.8125
$ws.Range("I126").Value = 20318.5
$ws.Range("K126").Value = 60955.5
$ws.Range("M126").Value = -58485.5
$ws.Range("N8").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3068.5
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 90000
$ws.Range("N132").Value = -95060

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 2500
$ws.Range("I41").Value = 2500
$ws.Range("K41").Value = 2500
$ws.Range("M41").Value = -2145
$ws.Range("H58").Value = 26666.334
$ws.Range("I58").Value = 24999.5
$ws.Range("K58").Value = 24999.5
$ws.Range("M58").Value = -24722.5
$ws.Range("H102").Value = 6670.6294
$ws.Range("I102").Value = 6255.3184
$ws.Range("K102").Value = 6255.3184
$ws.Range("M102").Value = -4633.3184
$ws.Range("H113").Value = 14173.75
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("H122").Value = 1002716.44
$ws.Range("I122").Value = 1836552.1
$ws.Range("K122").Value = 5509656.300000001
$ws.Range("M122").Value = -5507206.300000001
$ws.Range("H126").Value = 3860.2068
$ws.Range("J126").Value = 4695.4736
$ws.Range("L126").Value = 14086.4208
$ws.Range("N126").Value = -19026.4208
$ws.Range("H132").Value = 4514.4287
$ws.Range("I132").Value = 3528.8
$ws.Range("K132").Value = 10586.4
$ws.Range("M132").Value = -8056.400000000001
$ws.Range("M113").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5167.9614
$ws.Range("I7").Value = 3918.139
$ws.Range("K7").Value = 3918.139
$ws.Range("M7").Value = -3806.139
$ws.Range("H40").Value = 5629.3076
$ws.Range("I40").Value = 4658.125
$ws.Range("K40").Value = 4658.125
$ws.Range("M40").Value = -4522.125
$ws.Range("H61").Value = 3809.818
$ws.Range("I61").Value = 4338.625
$ws.Range("J61").Value = 2399.6667
$ws.Range("K61").Value = 4338.625
$ws.Range("L61").Value = 2399.6667
$ws.Range("M61").Value = -4136.625
$ws.Range("N61").Value = -2803.6667
$ws.Range("H93").Value = 956.4706
$ws.Range("I93").Value = 861.4286
$ws.Range("J93").Value = 1400
$ws.Range("K93").Value = 861.4286
$ws.Range("L93").Value = 1400
$ws.Range("M93").Value = 386.5714
$ws.Range("N93").Value = -3896
$ws.Range("H113").Value = 3809.818
$ws.Range("I113").Value = 4338.625
$ws.Range("J113").Value = 2399.6667
$ws.Range("K113").Value = 4338.625
$ws.Range("L113").Value = 2399.6667
$ws.Range("M113").Value = -2168.625
$ws.Range("N113").Value = -6739.6667
$ws.Range("H122").Value = 76197100
$ws.Range("I122").Value = 142864210
$ws.Range("K122").Value = 428592630
$ws.Range("M122").Value = -428590180
$ws.Range("H126").Value = 5167.9614
$ws.Range("I126").Value = 3918.139
$ws.Range("K126").Value = 11754.417
$ws.Range("M126").Value = -9284.417000000001
$ws.Range("H136").Value = 7706.2
$ws.Range("I136").Value = 10000
$ws.Range("K136").Value = 30000
$ws.Range("M136").Value = -27450

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28500
$ws.Range("J41").Value = 28500
$ws.Range("L41").Value = 28500
$ws.Range("N41").Value = -29280
$ws.Range("H126").Value = 2060.2
$ws.Range("I126").Value = 1933.6666
$ws.Range("K126").Value = 5800.9998
$ws.Range("M126").Value = -3330.9998
$ws.Range("H132").Value = 20838432
$ws.Range("I132").Value = 2069.3
$ws.Range("K132").Value = 6207.900000000001
$ws.Range("M132").Value = -3677.900000000001
